$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row of data for "2021 - Høst" exam right after the "2021 - Vår" row
$ws.Range("A28").Value = "2021 - Høst"
$ws.Range("B28").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-21-h.pdf)"
$ws.Range("C28").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-21-h-fasit.pdf)"

# Move selection down to mirror the recorded cursor position after the edit
$ws.Range("C29").Select()
